$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 48, shifting rows 48:59 down to 49:60
$ws.Rows.Item(48).Insert()

# Populate the newly inserted row 48 with the new weekly price record
$ws.Range("A48").Value = 9
$ws.Range("B48").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C48").Value = "Metropolitana"
$ws.Range("D48").Value = 44510
$ws.Range("E48").Value = 13
$ws.Range("F48").Value = "Fruta"
$ws.Range("G48").Value = 100102
$ws.Range("H48").Value = "Cítricos"
$ws.Range("I48").Value = 100102006
$ws.Range("J48").Value = "Pomelo"
$ws.Range("K48").Value = "Start Ruby"
$ws.Range("L48").Value = "Primera"
$ws.Range("M48").Value = 350
$ws.Range("N48").Value = 8000
$ws.Range("O48").Value = 8000
$ws.Range("P48").Value = 8000
$ws.Range("Q48").Value = "$/caja 14 kilos granel"
$ws.Range("R48").Value = "Región Metropolitana"
$ws.Range("S48").Value = 571
$ws.Range("T48").Value = 14

# Copy the date style from the row below (preserves the custom date/time number format)
$ws.Range("D49").Copy()
$ws.Range("D48").PasteSpecial(-4122)
